# The commit swaps the two embedded OOXML theme parts:
#   ppt/theme/theme1.xml (used by the slide master, "Integral")
#     -> becomes the stock "Office Theme" colour scheme
#   ppt/theme/theme2.xml (used only by the notes master, "Office Theme")
#     -> becomes the old "Integral" colour scheme
#
# i.e. the presentation's (slide-master) theme is switched from the
# "Integral" palette to the default "Office Theme" palette.  Re-create
# that by pushing the 12 "Office" theme colours into the presentation's
# ThemeColorScheme via the Design/Theme object model (COM surfaces the
# deck's single master theme through Presentation.SlideMaster.Theme).

function Convert-HexToOleRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Official "Office Theme" colour scheme (the values written into the
# post-commit theme1.xml <a:clrScheme>).
$officeColors = @(
    "000000", # 1 dk1
    "FFFFFF", # 2 lt1
    "44546A", # 3 dk2
    "E7E6E6", # 4 lt2
    "5B9BD5", # 5 accent1
    "ED7D31", # 6 accent2
    "A5A5A5", # 7 accent3
    "FFC000", # 8 accent4
    "4472C4", # 9 accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $colorScheme.Item($i).RGB = Convert-HexToOleRgb $officeColors[$i - 1]
}
